$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), copying the header style (s="1") from E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# time_taken values for data rows 2-92
$timeTaken = @(
    "2021-10-05 13:41:19.033036",
    "2021-10-05 13:41:19.033047",
    "2021-10-05 13:41:19.033050",
    "2021-10-05 13:41:19.033053",
    "2021-10-05 13:41:19.033056",
    "2021-10-05 13:41:19.033058",
    "2021-10-05 13:41:19.033061",
    "2021-10-05 13:41:19.033063",
    "2021-10-05 13:41:19.033066",
    "2021-10-05 13:41:19.033068",
    "2021-10-05 13:41:19.033071",
    "2021-10-05 13:41:19.033073",
    "2021-10-05 13:41:19.033076",
    "2021-10-05 13:41:19.033078",
    "2021-10-05 13:41:19.033081",
    "2021-10-05 13:41:19.033083",
    "2021-10-05 13:41:19.033086",
    "2021-10-05 13:41:19.033089",
    "2021-10-05 13:41:19.033091",
    "2021-10-05 13:41:19.033094",
    "2021-10-05 13:41:19.033096",
    "2021-10-05 13:41:19.033099",
    "2021-10-05 13:41:19.033102",
    "2021-10-05 13:41:19.033104",
    "2021-10-05 13:41:19.033107",
    "2021-10-05 13:41:19.033110",
    "2021-10-05 13:41:19.033112",
    "2021-10-05 13:41:19.033115",
    "2021-10-05 13:41:19.033117",
    "2021-10-05 13:41:19.033120",
    "2021-10-05 13:41:19.033122",
    "2021-10-05 13:41:19.033125",
    "2021-10-05 13:41:19.033128",
    "2021-10-05 13:41:19.033131",
    "2021-10-05 13:41:19.033133",
    "2021-10-05 13:41:19.033135",
    "2021-10-05 13:41:19.033138",
    "2021-10-05 13:41:19.033140",
    "2021-10-05 13:41:19.033143",
    "2021-10-05 13:41:19.033145",
    "2021-10-05 13:41:19.033148",
    "2021-10-05 13:41:19.033151",
    "2021-10-05 13:41:19.033153",
    "2021-10-05 13:41:19.033156",
    "2021-10-05 13:41:19.033158",
    "2021-10-05 13:41:19.033161",
    "2021-10-05 13:41:19.033163",
    "2021-10-05 13:41:19.033166",
    "2021-10-05 13:41:19.033168",
    "2021-10-05 13:41:19.033170",
    "2021-10-05 13:41:19.033173",
    "2021-10-05 13:41:19.033176",
    "2021-10-05 13:41:19.033179",
    "2021-10-05 13:41:19.033181",
    "2021-10-05 13:41:19.033184",
    "2021-10-05 13:41:19.033186",
    "2021-10-05 13:41:19.033189",
    "2021-10-05 13:41:19.033191",
    "2021-10-05 13:41:19.033194",
    "2021-10-05 13:41:19.033196",
    "2021-10-05 13:41:19.033199",
    "2021-10-05 13:41:19.033201",
    "2021-10-05 13:41:19.033204",
    "2021-10-05 13:41:19.033206",
    "2021-10-05 13:41:19.033210",
    "2021-10-05 13:41:19.033213",
    "2021-10-05 13:41:19.033216",
    "2021-10-05 13:41:19.033218",
    "2021-10-05 13:41:19.033220",
    "2021-10-05 13:41:19.033223",
    "2021-10-05 13:41:19.033225",
    "2021-10-05 13:41:19.033228",
    "2021-10-05 13:41:19.033231",
    "2021-10-05 13:41:19.033233",
    "2021-10-05 13:41:19.033236",
    "2021-10-05 13:41:19.033238",
    "2021-10-05 13:41:19.033243",
    "2021-10-05 13:41:19.033246",
    "2021-10-05 13:41:19.033248",
    "2021-10-05 13:41:19.033251",
    "2021-10-05 13:41:19.033253",
    "2021-10-05 13:41:19.033256",
    "2021-10-05 13:41:19.033258",
    "2021-10-05 13:41:19.033261",
    "2021-10-05 13:41:19.033263",
    "2021-10-05 13:41:19.033266",
    "2021-10-05 13:41:19.033269",
    "2021-10-05 13:41:19.033271",
    "2021-10-05 13:41:19.033274",
    "2021-10-05 13:41:19.033277",
    "2021-10-05 13:41:19.033279"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
